$d = $word.ActiveDocument

# The only substantive (visible) edit in the target revision is a single
# word swap: "rudimentary" -> "tentative" in the scheduling paragraph
# ("We have also set up a rudimentary schedule to track our progress
# throughout this project, ..."). Everything else in the supplied diff
# (the Jupyter / .ipynb / "Kaggle file format" / Github spans) is a no-op
# at the text level: those hunks only show Word's proofing engine merging
# runs that used to be split apart by <w:proofErr> spell-check markers --
# the rendered text before and after is identical there, so nothing needs
# to change in those paragraphs.
$found = $d.Content.Find.Execute("rudimentary", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "tentative", 2)

Write-Output "Replaced 'rudimentary' -> 'tentative': $found"
